$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate the "K" column (G) values — replaces the old "Strike#" values.
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 1
$ws.Range("G7").Value = 5
$ws.Range("G8").Value = 3
$ws.Range("G10").Value = 2
